$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.823.12'
$ws.Range("E2").Value = '  -1.82%  '
$ws.Range("D3").Value = '3.298.29'
$ws.Range("E3").Value = '  -2.05%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.582'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.88%  '
$ws.Range("D9").Value = '3.294.03'
$ws.Range("E9").Value = '  -1.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.184'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.21%  '
$ws.Range("E11").Value = '  -2.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.32'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.94%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000268'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.38%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.829.28'
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("B16").Value = 'BitcoinCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '626.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.41%  '
$ws.Range("D18").Value = '65.862.05'
$ws.Range("E18").Value = '  -1.89%  '
$ws.Range("E19").Value = '  -1.81%  '
$ws.Range("D20").Value = '3.304.91'
$ws.Range("E20").Value = '  -1.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.904'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '101.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  -1.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.51'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.62'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '546.70'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.36%  '
$ws.Range("E36").Value = '  -0.80%  '
$ws.Range("D37").Value = '3.804.81'
$ws.Range("E37").Value = '  +0.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '57.47'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.56%  '
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("D40").Value = '0.0₃0734'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '33.56'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.128'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("E43").Value = '  -6.68%  '
$ws.Range("E44").Value = '  -2.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.333'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.53%  '
$ws.Range("B46").Value = 'CoreDAO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -15.89%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0416'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.53%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.27'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.128'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.08%  '
